# "first of many passes through in cleaning database"
# Column I (rows 2-37) currently holds the text "No" for every sample row.
# Convert it to a real boolean (FALSE) value, formatted to still read as
# TRUE/FALSE text, instead of the free-text "No" string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 37; $r++) {
    $cell = $ws.Range("I$r")
    $cell.Value = $false
    $cell.NumberFormat = '"TRUE";"TRUE";"FALSE"'
}

# Move the active selection from column H to column I, matching where the
# editor's cursor ended up after cleaning that column.
$ws.Range("I2:I37").Select()

# Scroll the view down so the bottom of the table is visible.
$win = $excel.ActiveWindow
$win.ScrollRow = 33
$win.ScrollColumn = 1
